# Generate Report for handback
#
# For each locale sheet (zh-cn, de-de):
#   - Status ("B" column, shared across rows) flips from
#     "Ready for handoff" -> "Handed back: in sync with en-US"
#   - A new pair of columns is populated for each tracked source file:
#       "Latest Target File"   (E) = same display text/link as the source .md (A)
#       "Latest Handback File" (F) = same display text/link as the latest handoff .xlf (C)
#   - "Latest Handback DateTime" (G) is stamped with the handback time
#
# Row 2 -> 2b2c6534-1212-48d4-bcda-c18c04c8cfab.md
# Row 3 -> ad3a3400-23b0-4e83-a0f0-7e7181ce3337.md

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

function Set-HandbackColumns {
    param($ws, $row, $mdDisplay, $mdUrl, $xlfDisplay, $xlfUrl)

    $eCell = $ws.Range("E" + $row)
    $ws.Hyperlinks.Add($eCell, $mdUrl, "", "", $mdDisplay)
    $eCell.Style = "HyperLink"

    $fCell = $ws.Range("F" + $row)
    $ws.Hyperlinks.Add($fCell, $xlfUrl, "", "", $xlfDisplay)
    $fCell.Style = "HyperLink"
}

# ---------------- zh-cn ----------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B2").Value = $statusHandedBack
$wsZh.Range("B3").Value = $statusHandedBack

Set-HandbackColumns $wsZh "2" `
    "2b2c6534-1212-48d4-bcda-c18c04c8cfab.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/a5a14353733ec024090cd0b6ad854cc5b469e0db/e2e/2b2c6534-1212-48d4-bcda-c18c04c8cfab.md" `
    "2b2c6534-1212-48d4-bcda-c18c04c8cfab.c79eefe955e9552b6774a7a9738b25e785a78807.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2ce260ab52774affb21b7eda133dbe9c30fb98a3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/2b2c6534-1212-48d4-bcda-c18c04c8cfab.c79eefe955e9552b6774a7a9738b25e785a78807.zh-cn.xlf"

$wsZh.Range("G2").Value = "2016-02-15 08:16:21"

Set-HandbackColumns $wsZh "3" `
    "ad3a3400-23b0-4e83-a0f0-7e7181ce3337.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/a5a14353733ec024090cd0b6ad854cc5b469e0db/e2e/ad3a3400-23b0-4e83-a0f0-7e7181ce3337.md" `
    "ad3a3400-23b0-4e83-a0f0-7e7181ce3337.deccf20e52a181b41c40a9c995d5bb5f1d10971e.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2ce260ab52774affb21b7eda133dbe9c30fb98a3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/ad3a3400-23b0-4e83-a0f0-7e7181ce3337.deccf20e52a181b41c40a9c995d5bb5f1d10971e.zh-cn.xlf"

$wsZh.Range("G3").Value = "2016-02-15 08:16:21"

# ---------------- de-de ----------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B2").Value = $statusHandedBack
$wsDe.Range("B3").Value = $statusHandedBack

Set-HandbackColumns $wsDe "2" `
    "2b2c6534-1212-48d4-bcda-c18c04c8cfab.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/a5a14353733ec024090cd0b6ad854cc5b469e0db/e2e/2b2c6534-1212-48d4-bcda-c18c04c8cfab.md" `
    "2b2c6534-1212-48d4-bcda-c18c04c8cfab.c79eefe955e9552b6774a7a9738b25e785a78807.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e9f5a01c68544d466aeab0af0bd440fae3fd86c3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/2b2c6534-1212-48d4-bcda-c18c04c8cfab.c79eefe955e9552b6774a7a9738b25e785a78807.de-de.xlf"

$wsDe.Range("G2").Value = "2016-02-15 08:16:51"

Set-HandbackColumns $wsDe "3" `
    "ad3a3400-23b0-4e83-a0f0-7e7181ce3337.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/a5a14353733ec024090cd0b6ad854cc5b469e0db/e2e/ad3a3400-23b0-4e83-a0f0-7e7181ce3337.md" `
    "ad3a3400-23b0-4e83-a0f0-7e7181ce3337.deccf20e52a181b41c40a9c995d5bb5f1d10971e.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e9f5a01c68544d466aeab0af0bd440fae3fd86c3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/ad3a3400-23b0-4e83-a0f0-7e7181ce3337.deccf20e52a181b41c40a9c995d5bb5f1d10971e.de-de.xlf"

$wsDe.Range("G3").Value = "2016-02-15 08:16:51"
